$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "interview_type" column between username and test_taken ---
$ws.Columns.Item(2).Insert()

# --- Insert "final_rating" column between test_taken and answer_1 ---
$ws.Columns.Item(4).Insert()

# Header row
$ws.Range("B1").Value = "interview_type"
$ws.Range("D1").Value = "final_rating"

# Apply header style (bold/border/center) to the two new header cells, matching existing header style
$ws.Range("A1").Copy()
$ws.Range("B1,D1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Update data rows: usernames + interview_type values ---
$ws.Range("A2").Value = "user1"
$ws.Range("B2").Value = "Static"

$ws.Range("A3").Value = "user2"
$ws.Range("B3").Value = "Dynamic"

$ws.Range("A4").Value = "user3"
$ws.Range("B4").Value = "Hybrid"

# --- New row 7, D7 = single space ---
$ws.Range("D7").Value = " "

# --- Selection ---
$ws.Range("G14").Select()
